$wb = $excel.ActiveWorkbook

# Rename the "Include from mCSD Type codes " sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from mCSD Type codes ")
$includeSheet.Name = "Include #0"

$ws = $wb.Worksheets.Item("Metadata")

# Update Version: 3.8.0 -> 3.9.0
$ws.Range("B3").Value = "3.9.0"

# Update Experimental: (blank) -> "false" (must be stored as text, not boolean)
$helper = $ws.Range("Z100")
$helper.Formula = "=""false"""
$helper.Copy()
$ws.Range("B7").PasteSpecial(-4163, -4142, $false, $false)
$helper.Clear()

# Update Date
$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Update Contact rows
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Update Jurisdiction: World -> Global (Whole world)
$ws.Range("B13").Value = "Global (Whole world)"
